$data = @(
    @("hip",0,0,-0.5),
    @("hip",0,0.1,-0.5),
    @("hip",0,0.3,-0.5),
    @("hip",0,1,-0.5),
    @("hip",0,2,-0.3),
    @("hip",0,4,-0.2),
    @("hip",0,1000,-0.2),
    @("hip",9.9999,0,-0.5),
    @("hip",9.9999,0.1,-0.5),
    @("hip",9.9999,0.3,-0.5),
    @("hip",9.9999,1,-0.5),
    @("hip",9.9999,2,-0.3),
    @("hip",9.9999,4,-0.2),
    @("hip",9.9999,1000,-0.2),
    @("hip",10,0,-0.3),
    @("hip",10,0.1,-0.3),
    @("hip",10,0.3,-0.3),
    @("hip",10,1,-0.3),
    @("hip",10,2,-0.3),
    @("hip",10,4,-0.3),
    @("hip",10,1000,-0.3),
    @("hip",15,0,-0.3),
    @("hip",15,0.1,-0.3),
    @("hip",15,0.3,-0.3),
    @("hip",15,1,-0.3),
    @("hip",15,2,-0.3),
    @("hip",15,4,-0.3),
    @("hip",15,1000,-0.3),
    @("hip",20,0,-0.4),
    @("hip",20,0.1,-0.4),
    @("hip",20,0.3,-0.4),
    @("hip",20,1,-0.4),
    @("hip",20,2,-0.4),
    @("hip",20,4,-0.4),
    @("hip",20,1000,-0.4),
    @("hip",25,0,-0.75),
    @("hip",25,0.1,-0.75),
    @("hip",25,0.3,-0.5),
    @("hip",25,1,-0.5),
    @("hip",25,2,-0.5),
    @("hip",25,4,-0.5),
    @("hip",25,1000,-0.5),
    @("hip",90,0,-0.75),
    @("hip",90,0.1,-0.75),
    @("hip",90,0.3,-0.5),
    @("hip",90,1,-0.5),
    @("hip",90,2,-0.5),
    @("hip",90,4,-0.5),
    @("hip",90,1000,-0.5),
    @("gable",0,0,-0.5),
    @("gable",0,0.1,-0.5),
    @("gable",0,0.3,-0.5),
    @("gable",0,1,-0.5),
    @("gable",0,2,-0.3),
    @("gable",0,4,-0.2),
    @("gable",0,1000,-0.2),
    @("gable",90,0,-0.5),
    @("gable",90,0.1,-0.5),
    @("gable",90,0.3,-0.5),
    @("gable",90,1,-0.5),
    @("gable",90,2,-0.3),
    @("gable",90,4,-0.2),
    @("gable",90,1000,-0.2)
)

$wb = $excel.ActiveWorkbook

# Insert the new "cpe_t5_2b" worksheet right after "k_a" (before "cpe_t5_2c")
$wsAfter = $wb.Worksheets.Item("k_a")
$new = $wb.Worksheets.Add($null, $wsAfter)
$new.Name = "cpe_t5_2b"

$n = $data.Count

# Populate column by column so new shared strings are registered in the
# same order as the authored workbook (hip_or_gable, hip, gable, roof_pitch).

# Column B: hip_or_gable
$new.Cells.Item(1, 2).Value = "hip_or_gable"
for ($i = 0; $i -lt $n; $i++) {
    $new.Cells.Item($i + 2, 2).Value = $data[$i][0]
}

# Column C: roof_pitch
$new.Cells.Item(1, 3).Value = "roof_pitch"
for ($i = 0; $i -lt $n; $i++) {
    $new.Cells.Item($i + 2, 3).Value = $data[$i][1]
}

# Column D: d_b_ratio
$new.Cells.Item(1, 4).Value = "d_b_ratio"
for ($i = 0; $i -lt $n; $i++) {
    $new.Cells.Item($i + 2, 4).Value = $data[$i][2]
}

# Column E: c_pe
$new.Cells.Item(1, 5).Value = "c_pe"
for ($i = 0; $i -lt $n; $i++) {
    $new.Cells.Item($i + 2, 5).Value = $data[$i][3]
}

# Column A: version
$new.Cells.Item(1, 1).Value = "version"
for ($i = 0; $i -lt $n; $i++) {
    $new.Cells.Item($i + 2, 1).Value = "AS/NZS1170.2-2021"
}

# Column widths (best-fit on the lookup columns, fixed width on the rest)
$new.Columns.Item(1).ColumnWidth = 18.140625
$new.Columns.Item(2).ColumnWidth = 12.5703125
$new.Columns.Item(3).ColumnWidth = 14.7109375
$new.Columns.Item(4).ColumnWidth = 14.7109375
$new.Columns.Item(5).ColumnWidth = 11.140625

# Scroll/selection state matching the authored view
$new.Application.ActiveWindow.ScrollRow = 34
$new.Range("D58").Select()

Write-Output "done"
